$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:K1)
$headers = @("venue", "date", "result", "ownTeam", "oppTeam", "batsman", "totalRuns", "totalBalls", "total4s", "total6s", "sr")

# Data rows (A2:K7) - order/content per the updated data set
$data = @(
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Ravichandran Ashwin ", "7", "5", "1", "0", "140.00"),
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Delhi Capitals", "Rajasthan Royals", "Ravichandran Ashwin ", "0", "0", "0", "0", "-"),
    @(" Dubai (DSC)", " October 31 2020", "Mumbai won by 9 wickets (with 34 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Ravichandran Ashwin ", "12", "9", "0", "1", "133.33"),
    @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Delhi Capitals", "Rajasthan Royals", "Ravichandran Ashwin ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Delhi Capitals", "Kings XI Punjab", "Ravichandran Ashwin ", "4", "6", "0", "0", "66.66"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Delhi Capitals", "Kolkata Knight Riders", "Ravichandran Ashwin ", "14", "13", "2", "0", "107.69")
)

# The full used range is now A1:K7 - force it to be stored as text so
# numeric-looking values (e.g. "7", "140.00", "-") stay text, matching
# the original t="str" cell typing.
$fullRange = $ws.Range("A1:K7")
$fullRange.NumberFormat = "@"

# Write header row
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Write data rows
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
